$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Cells.Item(2, 4).Value = 0.0131
$ws.Cells.Item(5, 4).Value = -0.0798
$ws.Cells.Item(5, 5).Value = -0.0817
$ws.Cells.Item(5, 6).Value = -0.3785
$ws.Cells.Item(6, 4).Value = -0.0259
$ws.Cells.Item(6, 13).Value = -0.0442
$ws.Cells.Item(7, 4).Value = -0.124
$ws.Cells.Item(7, 5).Value = -0.1277
$ws.Cells.Item(7, 6).Value = 0.0537
$ws.Cells.Item(7, 7).Value = 0.1487
$ws.Cells.Item(8, 4).Value = -0.2996
$ws.Cells.Item(8, 5).Value = -0.3096
$ws.Cells.Item(8, 6).Value = -0.3441
$ws.Cells.Item(8, 7).Value = -0.6561
$ws.Cells.Item(8, 8).Value = -0.5964
$ws.Cells.Item(8, 9).Value = -0.378
$ws.Cells.Item(8, 10).Value = -0.1506
$ws.Cells.Item(8, 11).Value = -0.0479
$ws.Cells.Item(8, 12).Value = 0.0201
$ws.Cells.Item(9, 5).Value = 0.0661
$ws.Cells.Item(9, 6).Value = -0.0202
$ws.Cells.Item(9, 7).Value = -0.0066
$ws.Cells.Item(9, 8).Value = -0.151
$ws.Cells.Item(9, 9).Value = -0.0763
$ws.Cells.Item(9, 10).Value = -0.0441
$ws.Cells.Item(9, 11).Value = -0.0498
$ws.Cells.Item(9, 12).Value = -0.0215
$ws.Cells.Item(9, 13).Value = -0.0316
$ws.Cells.Item(10, 5).Value = 0.2461
$ws.Cells.Item(10, 6).Value = 0.212
$ws.Cells.Item(11, 4).Value = 0.0123
$ws.Cells.Item(12, 4).Value = 0.0119
$ws.Cells.Item(12, 5).Value = 0.5091
$ws.Cells.Item(12, 6).Value = -0.7283
$ws.Cells.Item(12, 7).Value = 0.0792
$ws.Cells.Item(12, 8).Value = -0.7729
$ws.Cells.Item(12, 9).Value = -0.5649
$ws.Cells.Item(13, 5).Value = -0.1017
$ws.Cells.Item(13, 6).Value = -0.1144
$ws.Cells.Item(13, 7).Value = -0.1294
$ws.Cells.Item(13, 8).Value = -0.1183
$ws.Cells.Item(13, 9).Value = -0.0503
$ws.Cells.Item(13, 10).Value = -0.0437
$ws.Cells.Item(13, 11).Value = -0.6123
$ws.Cells.Item(13, 12).Value = -0.6469
$ws.Cells.Item(13, 13).Value = -0.3046
$ws.Cells.Item(14, 4).Value = -0.9252
$ws.Cells.Item(14, 5).Value = -1.1031
$ws.Cells.Item(14, 6).Value = -0.9958
$ws.Cells.Item(14, 7).Value = -1.485
$ws.Cells.Item(14, 11).Value = -0.3158
$ws.Cells.Item(15, 4).Value = -2.3217
$ws.Cells.Item(15, 5).Value = -1.6943
$ws.Cells.Item(15, 6).Value = -1.8474
$ws.Cells.Item(15, 7).Value = -3.2671
$ws.Cells.Item(15, 8).Value = -2.8851
$ws.Cells.Item(15, 9).Value = -2.1119
$ws.Cells.Item(15, 10).Value = -1.9074
$ws.Cells.Item(15, 11).Value = -3.0357
$ws.Cells.Item(15, 12).Value = -1.2197
$ws.Cells.Item(15, 13).Value = -0.6379
$ws.Cells.Item(16, 4).Value = 23202.3
$ws.Cells.Item(17, 4).Value = 0.4181
$ws.Cells.Item(17, 5).Value = -0.0724
$ws.Cells.Item(17, 6).Value = -0.3114
$ws.Cells.Item(19, 5).Value = -1.1619
$ws.Cells.Item(19, 6).Value = -0.0021
$ws.Cells.Item(19, 8).Value = -0.0395
$ws.Cells.Item(19, 9).Value = -0.0349
$ws.Cells.Item(19, 10).Value = -0.0225
$ws.Cells.Item(19, 11).Value = -1.0899
$ws.Cells.Item(19, 12).Value = -0.0003
$ws.Cells.Item(19, 13).Value = -0.0003
$ws.Cells.Item(20, 4).Value = -1.1108
$ws.Cells.Item(20, 5).Value = -0.3758
$ws.Cells.Item(20, 7).Value = -1.0284
$ws.Cells.Item(20, 9).Value = -0.0206
$ws.Cells.Item(21, 4).Value = -0.1059
$ws.Cells.Item(21, 5).Value = 0.5857
$ws.Cells.Item(21, 6).Value = 0.682
$ws.Cells.Item(22, 4).Value = -0.0141
$ws.Cells.Item(22, 6).Value = -0.0196
$ws.Cells.Item(22, 8).Value = -0.0199
$ws.Cells.Item(22, 9).Value = -0.0197
$ws.Cells.Item(23, 4).Value = 0.094
$ws.Cells.Item(23, 5).Value = 0.2513
$ws.Cells.Item(23, 6).Value = 0.0847
$ws.Cells.Item(23, 7).Value = 0.0065
$ws.Cells.Item(23, 8).Value = 0.0301
$ws.Cells.Item(24, 4).Value = -0.0787
$ws.Cells.Item(24, 5).Value = -0.0334
$ws.Cells.Item(24, 6).Value = -0.0371
$ws.Cells.Item(24, 7).Value = -0.1178
$ws.Cells.Item(24, 8).Value = -0.0738
$ws.Cells.Item(24, 9).Value = -0.0796
$ws.Cells.Item(24, 10).Value = -0.0514
$ws.Cells.Item(24, 11).Value = -0.0174
$ws.Cells.Item(24, 12).Value = -0.0098
$ws.Cells.Item(26, 4).Value = 158.2
$ws.Cells.Item(28, 4).Value = -0.2049
$ws.Cells.Item(28, 5).Value = -0.0706
$ws.Cells.Item(28, 11).Value = -0.0242
$ws.Cells.Item(29, 4).Value = 0.0002
$ws.Cells.Item(32, 4).Value = -0.015
$ws.Cells.Item(32, 5).Value = 0.1315
$ws.Cells.Item(32, 6).Value = -0.1298
$ws.Cells.Item(33, 4).Value = 0.0042
$ws.Cells.Item(34, 4).Value = -0.0008
$ws.Cells.Item(34, 5).Value = -0.0007
$ws.Cells.Item(34, 6).Value = -0.0323
$ws.Cells.Item(34, 7).Value = -0.0319
$ws.Cells.Item(35, 4).Value = 0.004
$ws.Cells.Item(35, 5).Value = 0.0041
$ws.Cells.Item(35, 6).Value = 0.002
$ws.Cells.Item(35, 7).Value = 0.002
$ws.Cells.Item(35, 8).Value = 0.002
$ws.Cells.Item(35, 9).Value = 0.0019
$ws.Cells.Item(35, 10).Value = 0.0018
$ws.Cells.Item(35, 11).Value = 0.0018
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(36, 5).Value = 0.0458
$ws.Cells.Item(36, 6).Value = -0.0132
$ws.Cells.Item(36, 7).Value = -0.0003
$ws.Cells.Item(36, 8).Value = -0.0158
$ws.Cells.Item(36, 9).Value = -0.0002
$ws.Cells.Item(36, 10).Value = -0.0002
$ws.Cells.Item(36, 11).Value = -0.0002
$ws.Cells.Item(36, 12).Value = -0.0002
$ws.Cells.Item(36, 13).Value = -0.0001
$ws.Cells.Item(37, 5).Value = -0.0002
$ws.Cells.Item(38, 4).Value = 0.0005
$ws.Cells.Item(39, 4).Value = -0.0023
$ws.Cells.Item(39, 5).Value = 0.0037
$ws.Cells.Item(39, 6).Value = -0.0011
$ws.Cells.Item(39, 7).Value = -0.001
$ws.Cells.Item(39, 9).Value = -0.0048
$ws.Cells.Item(40, 5).Value = 0.0188
$ws.Cells.Item(40, 6).Value = -0.0028
$ws.Cells.Item(40, 7).Value = -0.001
$ws.Cells.Item(40, 8).Value = -0.001
$ws.Cells.Item(40, 9).Value = -0.001
$ws.Cells.Item(40, 10).Value = -0.0002
$ws.Cells.Item(40, 11).Value = -0.0002
$ws.Cells.Item(40, 12).Value = -0.0002
$ws.Cells.Item(40, 13).Value = -0.001
$ws.Cells.Item(41, 4).Value = -0.0002
$ws.Cells.Item(41, 5).Value = 0.0005
$ws.Cells.Item(41, 6).Value = -0.0002
$ws.Cells.Item(41, 7).Value = -0.0002
$ws.Cells.Item(41, 8).Value = -0.0001
$ws.Cells.Item(41, 9).Value = -0.0001
$ws.Cells.Item(41, 10).Value = -0.0001
$ws.Cells.Item(42, 4).Value = -0.0313
$ws.Cells.Item(42, 5).Value = 0.0483
$ws.Cells.Item(42, 6).Value = -0.0571
$ws.Cells.Item(42, 7).Value = -0.0401
$ws.Cells.Item(42, 8).Value = -0.0154
$ws.Cells.Item(42, 9).Value = -0.0003
$ws.Cells.Item(42, 10).Value = 0.0063
$ws.Cells.Item(42, 11).Value = 0.0074
$ws.Cells.Item(42, 12).Value = 0.0049
$ws.Cells.Item(42, 13).Value = 0.0016
$ws.Cells.Item(43, 4).Value = 15.3
$ws.Cells.Item(44, 4).Value = 0.0001
$ws.Cells.Item(44, 5).Value = 0.1315
$ws.Cells.Item(44, 6).Value = -0.1298
$ws.Cells.Item(46, 5).Value = 0.0128
$ws.Cells.Item(46, 6).Value = -0.0035
$ws.Cells.Item(46, 7).Value = -0.0001
$ws.Cells.Item(46, 8).Value = -0.0042
$ws.Cells.Item(47, 4).Value = -0.0002
$ws.Cells.Item(47, 5).Value = 0.0002
$ws.Cells.Item(48, 4).Value = 0.0214
$ws.Cells.Item(48, 5).Value = -0.1316
$ws.Cells.Item(48, 6).Value = 0.13
$ws.Cells.Item(49, 4).Value = -0.0023
$ws.Cells.Item(50, 4).Value = -0.0449
$ws.Cells.Item(50, 5).Value = -0.0443
$ws.Cells.Item(50, 6).Value = -0.0117
$ws.Cells.Item(50, 7).Value = -0.0116
$ws.Cells.Item(50, 8).Value = 0.0282
$ws.Cells.Item(51, 4).Value = 0.0041
$ws.Cells.Item(51, 5).Value = 0.007
$ws.Cells.Item(51, 6).Value = 0.0049
$ws.Cells.Item(51, 7).Value = 0.0034
$ws.Cells.Item(51, 8).Value = 0.0033
$ws.Cells.Item(51, 9).Value = 0.0033
$ws.Cells.Item(51, 10).Value = 0.0033
$ws.Cells.Item(51, 11).Value = 0.0032
$ws.Cells.Item(51, 12).Value = 0.0014
$ws.Cells.Item(53, 4).Value = 0.2
